$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the display text of the two existing URL cells to the new IP
# (172.31.15.230 -> 172.31.23.107), keeping the existing hyperlink
# relationships (rId1 / rId2) in place.
$ws.Range("A2").Value = "http://172.31.23.107:8082/webapp/"
$ws.Range("A3").Value = "http://172.31.23.107:8083"

# Re-point the existing hyperlinks' targets to the new IP as well -
# updating .Value alone does not touch the hyperlink relationship target,
# so walk the existing Hyperlinks collection and fix each Address in place.
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Address()
    if ($addr -like "*172.31.15.230:8083*") {
        $h.Address = "http://172.31.23.107:8083"
    }
    elseif ($addr -like "*172.31.15.230:8082/webapp*") {
        $h.Address = "http://172.31.23.107:8082/webapp/"
    }
}

# Add the two new rows / URLs for the third host.
$ws.Range("A4").Value = "http://172.31.29.70:8082/webapp"
$ws.Range("A5").Value = "http://172.31.29.70:8083"

$ws.Hyperlinks.Add($ws.Range("A4"), "http://172.31.29.70:8082/webapp")
$ws.Range("A4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A5"), "http://172.31.29.70:8083")
$ws.Range("A5").Style = "Hyperlink"
